$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: WhatsApp "fee_receipt" payment-receipt template (Fast2SMS message_id 4587)
$ws.Range("A7").Value = 4587
$ws.Range("B7").Value = "fee_receipt"
$ws.Range("C7").Value = "chords music academy (+917981585309)"
$ws.Range("D7").Value = "UTILITY"
$ws.Range("E7").Value = "APPROVED"
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = "Var1|Var2|Var3|Var4|Var5|Var6"
$ws.Range("H7").Value = "No media required"
$ws.Range("I7").Value = "https://www.fast2sms.com/dev/whatsapp?authorization=<YOUR_API_KEY>&message_id=4587&numbers=<MOBILE_NUMBER>&variables_values=Var1|Var2|Var3|Var4|Var5|Var6"
$ws.Range("J7").Value = "BODY: Dear {Var1}, `nThank you for your payment to Chords Music Academy! 🎵`n💰 Payment Details:`n- Amount: ₹{Var2}`n- Receipt No: {Var3}`n- Package: {Var4}`n- Payment Date: {Var5}`n{Var6}`n🎶 Keep practicing and let your musical journey flourish!`n📞 Contact us at +91 7981585309 | BUTTON: Call (PHONE_NUMBER) - +917981585309"

# Emphasise the template name / variables_values columns to mirror the
# source workbook's highlight font used for these two cells.
$ws.Range("B7").Font.Bold = $false
$ws.Range("B7").Font.Italic = $false
$ws.Range("B7").Font.Strikethrough = $false
$ws.Range("B7").Font.Underline = -4142
$ws.Range("G7").Font.Bold = $false
$ws.Range("G7").Font.Italic = $false
$ws.Range("G7").Font.Strikethrough = $false
$ws.Range("G7").Font.Underline = -4142

# The long wrapped message body would otherwise force an auto-computed
# custom row height; restore the sheet's default row height like the
# source file.
$ws.Rows.Item(7).AutoFit()

$ws.Range("A7").Select()

Write-Host "done"
